# Automatic data update for establecimientos_pendientes.xlsx
# Applies the monthly "2024-11" (column O) data refresh: fills in newly
# reported monthly counts, recomputes the yearly "Total" (column D), and
# flips each establishment's "Estado" (column R) between "Al dia" and
# "Pendiente" to reflect which rows still have missing monthly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed / newly-reported values
$ws.Range("D2").Value = 15627
$ws.Range("O2").Value = 1401
$ws.Range("O3").Value = "Pendiente"
$ws.Range("R3").Value = "Pendiente"
$ws.Range("D4").Value = 38
$ws.Range("G4").Value = "Pendiente"
$ws.Range("K4").Value = 1
$ws.Range("P4").Value = 6
$ws.Range("R4").Value = "Pendiente"
$ws.Range("O5").Value = "Pendiente"
$ws.Range("F6").Value = "Pendiente"
$ws.Range("K6").Value = 28
$ws.Range("R6").Value = "Pendiente"
$ws.Range("O7").Value = "Pendiente"
$ws.Range("O9").Value = "Pendiente"
$ws.Range("D10").Value = 219
$ws.Range("P10").Value = 26
$ws.Range("O11").Value = "Pendiente"
$ws.Range("D12").Value = 3668
$ws.Range("O12").Value = 38
$ws.Range("D13").Value = 3640
$ws.Range("P13").Value = 224
$ws.Range("O14").Value = "Pendiente"
$ws.Range("R14").Value = "Pendiente"
$ws.Range("O15").Value = "Pendiente"
$ws.Range("D17").Value = 268
$ws.Range("P17").Value = 19
$ws.Range("O18").Value = "Pendiente"
$ws.Range("R18").Value = "Pendiente"
$ws.Range("D19").Value = 34325
$ws.Range("N19").Value = 3891
$ws.Range("O19").Value = "Pendiente"
$ws.Range("O21").Value = "Pendiente"
$ws.Range("D22").Value = 221
$ws.Range("P22").Value = 13
$ws.Range("D23").Value = 621
$ws.Range("P23").Value = 18
$ws.Range("O24").Value = "Pendiente"
$ws.Range("O27").Value = "Pendiente"
$ws.Range("R27").Value = "Pendiente"
$ws.Range("D28").Value = 153
$ws.Range("O28").Value = 15
$ws.Range("P28").Value = 8
$ws.Range("H29").Value = 2
$ws.Range("J29").Value = 1
$ws.Range("K29").Value = 2
$ws.Range("O29").Value = "Pendiente"
$ws.Range("R29").Value = "Pendiente"
$ws.Range("O30").Value = "Pendiente"
$ws.Range("O31").Value = "Pendiente"
$ws.Range("R31").Value = "Pendiente"
$ws.Range("D33").Value = 194
$ws.Range("M33").Value = "Pendiente"
$ws.Range("N33").Value = "Pendiente"
$ws.Range("P33").Value = 11
$ws.Range("R33").Value = "Pendiente"
$ws.Range("D34").Value = 936
$ws.Range("H34").Value = 183
$ws.Range("I34").Value = 120
$ws.Range("J34").Value = 155
$ws.Range("K34").Value = 4
$ws.Range("O34").Value = "Pendiente"
$ws.Range("F35").Value = "Pendiente"
$ws.Range("H35").Value = "Pendiente"
$ws.Range("K35").Value = 8
$ws.Range("R35").Value = "Pendiente"
$ws.Range("O36").Value = "Pendiente"
$ws.Range("R36").Value = "Pendiente"
$ws.Range("D37").Value = 29
$ws.Range("M37").Value = 2
$ws.Range("R37").Value = "Al día"
$ws.Range("D38").Value = 9278
$ws.Range("O38").Value = 777
$ws.Range("O39").Value = "Pendiente"
$ws.Range("O40").Value = "Pendiente"
$ws.Range("R40").Value = "Pendiente"
$ws.Range("O41").Value = "Pendiente"
$ws.Range("R41").Value = "Pendiente"
$ws.Range("O42").Value = "Pendiente"
$ws.Range("O45").Value = "Pendiente"
$ws.Range("D46").Value = 167
$ws.Range("L46").Value = 26
$ws.Range("O46").Value = "Pendiente"
$ws.Range("O47").Value = "Pendiente"
$ws.Range("R47").Value = "Pendiente"
$ws.Range("D49").Value = 1023
$ws.Range("P49").Value = 60
$ws.Range("D50").Value = 672
$ws.Range("N50").Value = 52
$ws.Range("O50").Value = 49
$ws.Range("R50").Value = "Al día"
$ws.Range("O51").Value = "Pendiente"
$ws.Range("R51").Value = "Pendiente"
$ws.Range("D52").Value = 12737
$ws.Range("N52").Value = 1629
$ws.Range("O52").Value = "Pendiente"
$ws.Range("F53").Value = "Pendiente"
$ws.Range("G53").Value = "Pendiente"
$ws.Range("H53").Value = "Pendiente"
$ws.Range("J53").Value = "Pendiente"
$ws.Range("M53").Value = "Pendiente"
$ws.Range("O53").Value = "Pendiente"
$ws.Range("F54").Value = "Pendiente"
$ws.Range("G54").Value = "Pendiente"
$ws.Range("H54").Value = "Pendiente"
$ws.Range("I54").Value = "Pendiente"
$ws.Range("J54").Value = "Pendiente"
$ws.Range("M54").Value = "Pendiente"
$ws.Range("O54").Value = "Pendiente"
$ws.Range("F55").Value = "Pendiente"
$ws.Range("G55").Value = "Pendiente"
$ws.Range("H55").Value = "Pendiente"
$ws.Range("O55").Value = "Pendiente"
$ws.Range("D56").Value = 145
$ws.Range("N56").Value = 15
$ws.Range("O56").Value = 5
$ws.Range("E57").Value = "Pendiente"
$ws.Range("F57").Value = "Pendiente"
$ws.Range("G57").Value = "Pendiente"
$ws.Range("H57").Value = "Pendiente"
$ws.Range("I57").Value = "Pendiente"
$ws.Range("J57").Value = "Pendiente"
$ws.Range("K57").Value = "Pendiente"
$ws.Range("L57").Value = "Pendiente"
$ws.Range("M57").Value = "Pendiente"
$ws.Range("N57").Value = "Pendiente"
$ws.Range("O57").Value = "Pendiente"
$ws.Range("R57").Value = "Pendiente"
$ws.Range("O59").Value = "Pendiente"
$ws.Range("R59").Value = "Pendiente"
$ws.Range("O60").Value = "Pendiente"
$ws.Range("R60").Value = "Pendiente"
$ws.Range("O61").Value = "Pendiente"
$ws.Range("R61").Value = "Pendiente"
$ws.Range("D62").Value = 3722
$ws.Range("N62").Value = 422
$ws.Range("O62").Value = "Pendiente"
$ws.Range("O63").Value = "Pendiente"
$ws.Range("D65").Value = 755
$ws.Range("L65").Value = 81
$ws.Range("M65").Value = 54
$ws.Range("N65").Value = 77
$ws.Range("O65").Value = 33
$ws.Range("R65").Value = "Al día"
$ws.Range("D66").Value = 225
$ws.Range("L66").Value = 36
$ws.Range("O66").Value = "Pendiente"
$ws.Range("O67").Value = "Pendiente"
$ws.Range("R67").Value = "Pendiente"
$ws.Range("O68").Value = "Pendiente"
$ws.Range("R68").Value = "Pendiente"
$ws.Range("E69").Value = "Pendiente"
$ws.Range("F69").Value = "Pendiente"
$ws.Range("G69").Value = "Pendiente"
$ws.Range("H69").Value = "Pendiente"
$ws.Range("I69").Value = "Pendiente"
$ws.Range("J69").Value = "Pendiente"
$ws.Range("K69").Value = "Pendiente"
$ws.Range("L69").Value = "Pendiente"
$ws.Range("M69").Value = "Pendiente"
$ws.Range("N69").Value = "Pendiente"
$ws.Range("O69").Value = "Pendiente"
$ws.Range("R69").Value = "Pendiente"
$ws.Range("O70").Value = "Pendiente"
$ws.Range("O71").Value = "Pendiente"
$ws.Range("E72").Value = "Pendiente"
$ws.Range("F72").Value = "Pendiente"
$ws.Range("G72").Value = "Pendiente"
$ws.Range("H72").Value = "Pendiente"
$ws.Range("I72").Value = "Pendiente"
$ws.Range("J72").Value = "Pendiente"
$ws.Range("K72").Value = "Pendiente"
$ws.Range("O72").Value = "Pendiente"

# These cells are fully removed (no longer present) in the updated sheet
$ws.Range("P53").ClearContents()
$ws.Range("P57").ClearContents()
$ws.Range("P69").ClearContents()

